# New crime data collected — refresh the weekly CompStat figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: report volume/number and the reporting week dates.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# ---------------------------------------------------------------------------
# Helper: write a bunch of numeric cell values for a row in one shot.
# ---------------------------------------------------------------------------
function Set-Cells {
    param($Sheet, [hashtable]$Values)
    foreach ($addr in $Values.Keys) {
        $Sheet.Range($addr).Value = $Values[$addr]
    }
}

# Row 16 — Robbery
Set-Cells $ws @{
    "F16" = 8
    "G16" = 1
    "H16" = 700
    "I16" = 54
    "K16" = 1.886792452830
    "L16" = -27.027027027027
    "M16" = -22.857142857142
    "N16" = -80.851063829787
}

# Row 17 — Fel. Assault
Set-Cells $ws @{
    "C17" = 2
    "D17" = 3
    "E17" = -33.333333333333
    "F17" = 10
    "G17" = 9
    "H17" = 11.111111111111
    "I17" = 80
    "J17" = 87
    "K17" = -8.045977011494
    "L17" = -8.045977011494
    "M17" = 73.913043478260
    "N17" = -27.272727272727
}

# Row 18 — Burglary
Set-Cells $ws @{
    "C18" = 8
    "D18" = 12
    "E18" = -33.333333333333
    "F18" = 25
    "G18" = 34
    "H18" = -26.470588235294
    "I18" = 261
    "J18" = 272
    "K18" = -4.044117647058
    "L18" = 10.126582278481
    "M18" = 23.113207547169
    "N18" = -70.374574347332
}

# Row 19 — Gr. Larceny
Set-Cells $ws @{
    "C19" = 10
    "E19" = -28.571428571428
    "F19" = 39
    "G19" = 48
    "H19" = -18.75
    "I19" = 397
    "J19" = 544
    "K19" = -27.022058823529
    "L19" = -31.076388888888
    "M19" = 25.632911392405
    "N19" = -19.635627530364
}

# Row 20 — G.L.A.
Set-Cells $ws @{
    "C20" = 3
    "D20" = 4
    "E20" = -25
    "F20" = 23
    "G20" = 13
    "H20" = 76.923076923076
    "I20" = 225
    "J20" = 164
    "K20" = 37.195121951219
    "L20" = 127.272727272727
    "M20" = 85.950413223140
    "N20" = -92.135616917161
}

# Row 21 — TOTAL
Set-Cells $ws @{
    "C21" = 25
    "D21" = 33
    "E21" = -24.242424242424
    "F21" = 105
    "G21" = 105
    "H21" = 0
    "I21" = 1026
    "J21" = 1131
    "K21" = -9.283819628647
    "L21" = -4.558139534883
    "M21" = 33.246753246753
    "N21" = -77.902218393280
}

# Row 24 — Petit Larceny
Set-Cells $ws @{
    "C24" = 5
    "D24" = 14
    "E24" = -64.285714285714
    "F24" = 51
    "G24" = 46
    "H24" = 10.869565217391
    "I24" = 486
    "J24" = 505
    "K24" = -3.762376237623
    "L24" = -27.678571428571
    "M24" = 26.233766233766
}

# Row 25 — Retail Theft
Set-Cells $ws @{
    "C25" = 3
    "F25" = 16
    "G25" = 9
    "H25" = 77.777777777777
    "I25" = 107
    "K25" = 1.904761904761
    "L25" = 4.901960784313
}
# D25 and E25 flip from numbers to the "no data" placeholder text used
# throughout this sheet ("0" / "***.*" shared strings with the s="13"
# right-aligned text style). Copying an existing placeholder cell of the
# same kind reproduces that exact style/shared-string combination.
$ws.Range("D14").Copy($ws.Range("D25"))
$ws.Range("E14").Copy($ws.Range("E25"))

# Row 26 — Misd. Assault
Set-Cells $ws @{
    "C26" = 4
    "E26" = 33.333333333333
    "G26" = 16
    "H26" = 50
    "I26" = 195
    "J26" = 195
    "K26" = 0
    "L26" = -9.302325581395
    "M26" = 28.289473684210
}

# Row 28 — Other Sex Crimes
# C28 and G28/H28 flip from numbers to the placeholder text, F28 stays numeric.
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("F28").Value = 1
$ws.Range("C14").Copy($ws.Range("G28"))
$ws.Range("E14").Copy($ws.Range("H28"))

# Row 33 — Traffic Fatalities
$ws.Range("F33").Value = 1
